$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume strings are kept as text (matches original formatting)
$textCells = @("D2","E2","D3","E3","E4","D5","E5","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","D20","E20","D21","E21","E22","D23","E23","D24","E24","D25","E25","E26","D27","E27","E28","E29","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","E37","E38","E39","D40","E40","E41","E42","E43","D44","E44","E45","D46","E46","E47","D48","E48","E49","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the latest crypto data refresh
$ws.Range("D2").Value = "34.062.75"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.777.52"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "224.85"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "31.69"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").Value = "0.290"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "0.0683"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "2.034.77"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "10.86"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.755.31"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "34.073.01"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "0.619"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "4.17"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "67.46"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "244.04"
$ws.Range("D20").Value = "0.0₃0783"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "4.08"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "161.03"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").Value = "16.17"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D31").Value = "0.0514"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "3.70"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "3.70"
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").Value = "1.437.38"
$ws.Range("E35").Value = "  +3.14%  "
$ws.Range("D36").Value = "0.652"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +4.63%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "80.07"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "13.55"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "6.03"
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "1.936.89"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  -6.77%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "103.68"
$ws.Range("E51").Value = "  -3.48%  "
